$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.00391733333333
$ws.Range("H2").Value = 105.011752
$ws.Range("I2").Value = 0.9591895364534718
$ws.Range("J2").Value = 0.9591895364534718
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1178306666666667
$ws.Range("N2").Value = 0.353492
$ws.Range("O2").Value = 0.0005211703885903252
$ws.Range("P2").Value = 0.0005211703885903251
$ws.Range("Q2").Value = 4.124534915331556
$ws.Range("R2").Value = 37.12081423798401
$ws.Range("S2").Value = 0.0004999011834452298
$ws.Range("T2").Value = 0.0004999011834452297

$ws.Range("G3").Value = 35.00391733333333
$ws.Range("H3").Value = 105.011752
$ws.Range("I3").Value = 0.9591895364534718
$ws.Range("J3").Value = 0.9591895364534718
$ws.Range("O3").Value = 0.9986266812609277
$ws.Range("P3").Value = 0.9986266812609277
$ws.Range("Q3").Value = 7903.11710798305
$ws.Range("R3").Value = 71128.05397184745
$ws.Range("S3").Value = 0.9578722634887382
$ws.Range("T3").Value = 0.9578722634887382

$ws.Range("G4").Value = 35.00391733333333
$ws.Range("H4").Value = 105.011752
$ws.Range("I4").Value = 0.9591895364534718
$ws.Range("J4").Value = 0.9591895364534718
$ws.Range("M4").Value = 0.192661
$ws.Range("N4").Value = 0.5779829999999999
$ws.Range("O4").Value = 0.0008521483504820529
$ws.Range("P4").Value = 0.0008521483504820528
$ws.Range("Q4").Value = 6.743889717357332
$ws.Range("R4").Value = 60.69500745621599
$ws.Range("S4").Value = 0.000817371781288471
$ws.Range("T4").Value = 0.0008173717812884709

$ws.Range("I5").Value = 0.0008369499257158872
$ws.Range("J5").Value = 0.0008369499257158872
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1178306666666667
$ws.Range("N5").Value = 0.353492
$ws.Range("O5").Value = 0.0005211703885903252
$ws.Range("P5").Value = 0.0005211703885903251
$ws.Range("Q5").Value = 0.003598902052
$ws.Range("R5").Value = 0.032390118468
$ws.Range("S5").Value = [double]"4.361935180159928e-07"
$ws.Range("T5").Value = [double]"4.361935180159927e-07"

$ws.Range("I6").Value = 0.0008369499257158872
$ws.Range("J6").Value = 0.0008369499257158872
$ws.Range("O6").Value = 0.9986266812609277
$ws.Range("P6").Value = 0.9986266812609277
$ws.Range("S6").Value = 0.0008358005266992365
$ws.Range("T6").Value = 0.0008358005266992365

$ws.Range("I7").Value = 0.0008369499257158872
$ws.Range("J7").Value = 0.0008369499257158872
$ws.Range("M7").Value = 0.192661
$ws.Range("N7").Value = 0.5779829999999999
$ws.Range("O7").Value = 0.0008521483504820529
$ws.Range("P7").Value = 0.0008521483504820528
$ws.Range("Q7").Value = 0.005884444922999999
$ws.Range("R7").Value = 0.05296000430699999
$ws.Range("S7").Value = [double]"7.132054986348701e-07"
$ws.Range("T7").Value = [double]"7.1320549863487e-07"

$ws.Range("G8").Value = 1.458762333333333
$ws.Range("H8").Value = 4.376287
$ws.Range("I8").Value = 0.03997351362081222
$ws.Range("J8").Value = 0.03997351362081222
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1178306666666667
$ws.Range("N8").Value = 0.353492
$ws.Range("O8").Value = 0.0005211703885903252
$ws.Range("P8").Value = 0.0005211703885903251
$ws.Range("Q8").Value = 0.1718869382448889
$ws.Range("R8").Value = 1.546982444204
$ws.Range("S8").Value = [double]"2.083301162707936e-05"
$ws.Range("T8").Value = [double]"2.083301162707936e-05"

$ws.Range("G9").Value = 1.458762333333333
$ws.Range("H9").Value = 4.376287
$ws.Range("I9").Value = 0.03997351362081222
$ws.Range("J9").Value = 0.03997351362081222
$ws.Range("O9").Value = 0.9986266812609277
$ws.Range("P9").Value = 0.9986266812609277
$ws.Range("Q9").Value = 329.3565529612706
$ws.Range("R9").Value = 2964.208976651435
$ws.Range("S9").Value = 0.03991861724549019
$ws.Range("T9").Value = 0.03991861724549019

$ws.Range("G10").Value = 1.458762333333333
$ws.Range("H10").Value = 4.376287
$ws.Range("I10").Value = 0.03997351362081222
$ws.Range("J10").Value = 0.03997351362081222
$ws.Range("M10").Value = 0.192661
$ws.Range("N10").Value = 0.5779829999999999
$ws.Range("O10").Value = 0.0008521483504820529
$ws.Range("P10").Value = 0.0008521483504820528
$ws.Range("Q10").Value = 0.2810466099023333
$ws.Range("R10").Value = 2.529419489121
$ws.Range("S10").Value = [double]"3.4063363694947e-05"
$ws.Range("T10").Value = [double]"3.4063363694947e-05"

Write-Host "Applied 100 cell updates"
